$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.041119126895842
$ws.Range("D2").Value = 1.045006794012703
$ws.Range("E2").Value = 1.049815510834324
$ws.Range("F2").Value = 1.061655576950616
$ws.Range("I2").Value = 1.044354288157395
$ws.Range("J2").Value = 1.046202055977866
$ws.Range("K2").Value = 1.047776245389759
$ws.Range("L2").Value = 1.052571518963095
$ws.Range("M2").Value = 1.064379047113664
$ws.Range("N2").Value = 1.047687782898408
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.042014604222189
$ws.Range("D3").Value = 1.045700229525413
$ws.Range("E3").Value = 1.05064641970567
$ws.Range("F3").Value = 1.062648406746973
$ws.Range("I3").Value = 1.044618088882875
$ws.Range("J3").Value = 1.04674370010143
$ws.Range("K3").Value = 1.048281341465999
$ws.Range("L3").Value = 1.053214701598445
$ws.Range("M3").Value = 1.06518609156722
$ws.Range("N3").Value = 1.048230196218755
$ws.Range("B4").Value = 1.019999999999999
$ws.Range("C4").Value = 1.042594450378797
$ws.Range("D4").Value = 1.046149186519731
$ws.Range("E4").Value = 1.051184847752058
$ws.Range("F4").Value = 1.063291880806394
$ws.Range("I4").Value = 1.044787678839789
$ws.Range("J4").Value = 1.047093943625774
$ws.Range("K4").Value = 1.048607753039079
$ws.Range("L4").Value = 1.053631008467366
$ws.Range("M4").Value = 1.065708731462802
$ws.Range("N4").Value = 1.048580937129077
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.04283831525401
$ws.Range("D5").Value = 1.046337988318112
$ws.Range("E5").Value = 1.051411386872007
$ws.Range("F5").Value = 1.063562646479921
$ws.Range("I5").Value = 1.044858708836618
$ws.Range("J5").Value = 1.047241128380312
$ws.Range("K5").Value = 1.048744875015406
$ws.Range("L5").Value = 1.053806052589218
$ws.Range("M5").Value = 1.065928550505637
$ws.Range("N5").Value = 1.048728330902849
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.042879266925543
$ws.Range("D6").Value = 1.046369692463338
$ws.Range("E6").Value = 1.051449434533425
$ws.Range("F6").Value = 1.063608123798535
$ws.Range("I6").Value = 1.044870619490269
$ws.Range("J6").Value = 1.047265837939615
$ws.Range("K6").Value = 1.048767892418989
$ws.Range("L6").Value = 1.053835444918259
$ws.Range("M6").Value = 1.065965464989421
$ws.Range("N6").Value = 1.04875307555256
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.042597708531368
$ws.Range("D7").Value = 1.046151709063573
$ws.Range("E7").Value = 1.05118787405779
$ws.Range("F7").Value = 1.063295497813253
$ws.Range("I7").Value = 1.044788628990617
$ws.Range("J7").Value = 1.047095910543361
$ws.Range("K7").Value = 1.048609585669581
$ws.Range("L7").Value = 1.053633347304087
$ws.Range("M7").Value = 1.065711668299437
$ws.Range("N7").Value = 1.048582906839913
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.041421671534567
$ws.Range("D8").Value = 1.045241089529176
$ws.Range("E8").Value = 1.050096159161011
$ws.Range("F8").Value = 1.061990890825652
$ws.Range("I8").Value = 1.0444436696094
$ws.Range("J8").Value = 1.046385155489856
$ws.Range("K8").Value = 1.047947031172494
$ws.Range("L8").Value = 1.052788859042454
$ws.Range("M8").Value = 1.064651702143762
$ws.Range("N8").Value = 1.047871142432706
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.039352555524132
$ws.Range("D9").Value = 1.043638500446225
$ws.Range("E9").Value = 1.048178412724757
$ws.Range("F9").Value = 1.059700086815532
$ws.Range("I9").Value = 1.04382735998144
$ws.Range("J9").Value = 1.045130948378695
$ws.Range("K9").Value = 1.046776360151693
$ws.Range("L9").Value = 1.051301768125253
$ws.Range("M9").Value = 1.062787241901278
$ws.Range("N9").Value = 1.046615154203589
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.037975373965457
$ws.Range("D10").Value = 1.042571564491027
$ws.Range("E10").Value = 1.046904030488268
$ws.Range("F10").Value = 1.058178400967619
$ws.Range("I10").Value = 1.043410846632912
$ws.Range("J10").Value = 1.044293683121915
$ws.Range("K10").Value = 1.045993846787538
$ws.Range("L10").Value = 1.05031111733691
$ws.Range("M10").Value = 1.061546586258634
$ws.Range("N10").Value = 1.045776699934108
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.037379582291667
$ws.Range("D11").Value = 1.042109932953209
$ws.Range("E11").Value = 1.046353202757381
$ws.Range("F11").Value = 1.057520820211498
$ws.Range("I11").Value = 1.043229162070341
$ws.Range("J11").Value = 1.043930883169907
$ws.Range("K11").Value = 1.045654532782139
$ws.Range("L11").Value = 1.049882344881816
$ws.Range("M11").Value = 1.061009934071966
$ws.Range("N11").Value = 1.045413384764554
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.037158360540783
$ws.Range("D12").Value = 1.041938517780339
$ws.Range("E12").Value = 1.046148750688746
$ws.Range("F12").Value = 1.057276765027582
$ws.Range("I12").Value = 1.043161476908284
$ws.Range("J12").Value = 1.043796085246923
$ws.Range("K12").Value = 1.045528425271391
$ws.Range("L12").Value = 1.049723108731034
$ws.Range("M12").Value = 1.060810682905496
$ws.Range("N12").Value = 1.045278395413058
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.037205809637428
$ws.Range("D13").Value = 1.041975284389806
$ws.Range("E13").Value = 1.046192599546601
$ws.Range("F13").Value = 1.057329106638978
$ws.Range("I13").Value = 1.043176004619747
$ws.Range("J13").Value = 1.043825001556803
$ws.Range("K13").Value = 1.045555478965294
$ws.Range("L13").Value = 1.049757264094491
$ws.Range("M13").Value = 1.060853419074993
$ws.Range("N13").Value = 1.045307352787415
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.037361294357439
$ws.Range("D14").Value = 1.042095762584457
$ws.Range("E14").Value = 1.046336299620009
$ws.Range("F14").Value = 1.057500642454317
$ws.Range("I14").Value = 1.043223571258356
$ws.Range("J14").Value = 1.043919741498347
$ws.Range("K14").Value = 1.045644110142152
$ws.Range("L14").Value = 1.049869181773318
$ws.Range("M14").Value = 1.060993462160136
$ws.Range("N14").Value = 1.045402227270543
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.037457104498587
$ws.Range("D15").Value = 1.042170000551658
$ws.Range("E15").Value = 1.046424857880247
$ws.Range("F15").Value = 1.057606357820307
$ws.Range("I15").Value = 1.04325285222706
$ws.Range("J15").Value = 1.043978108914724
$ws.Range("K15").Value = 1.045698709333641
$ws.Range("L15").Value = 1.049938141828838
$ws.Range("M15").Value = 1.061079758673523
$ws.Range("N15").Value = 1.045460677575345
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.038014926106904
$ws.Range("D16").Value = 1.042602209137098
$ws.Range("E16").Value = 1.046940608051162
$ws.Range("F16").Value = 1.058222070350033
$ws.Range("I16").Value = 1.043422876438109
$ws.Range("J16").Value = 1.044317755595729
$ws.Range("K16").Value = 1.046016355902605
$ws.Range("L16").Value = 1.050339577563036
$ws.Range("M16").Value = 1.061582213998296
$ws.Range("N16").Value = 1.045800806593596
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.038364977432999
$ws.Range("D17").Value = 1.042873419229277
$ws.Range("E17").Value = 1.047264390069977
$ws.Range("F17").Value = 1.05860864473194
$ws.Range("I17").Value = 1.043529172082609
$ws.Range("J17").Value = 1.044530738450243
$ws.Range("K17").Value = 1.046215479262556
$ws.Range("L17").Value = 1.050591438136405
$ws.Range("M17").Value = 1.061897541717497
$ws.Range("N17").Value = 1.046014091908192
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.038569208189192
$ws.Range("D18").Value = 1.043031646007742
$ws.Range("E18").Value = 1.047453341894582
$ws.Range("F18").Value = 1.05883425422483
$ws.Range("I18").Value = 1.043591044018108
$ws.Range("J18").Value = 1.044654942731774
$ws.Range("K18").Value = 1.046331578114709
$ws.Range("L18").Value = 1.050738362007652
$ws.Range("M18").Value = 1.06208152098002
$ws.Range("N18").Value = 1.04613847257405
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.0386388543412
$ws.Range("D19").Value = 1.043085603061862
$ws.Range("E19").Value = 1.047517785697083
$ws.Range("F19").Value = 1.058911202764257
$ws.Range("I19").Value = 1.043612118931841
$ws.Range("J19").Value = 1.044697288909418
$ws.Range("K19").Value = 1.046371156915959
$ws.Range("L19").Value = 1.050788462246661
$ws.Range("M19").Value = 1.062144262265596
$ws.Range("N19").Value = 1.046180878888123
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.038327414878716
$ws.Range("D20").Value = 1.042844317372145
$ws.Range("E20").Value = 1.047229641461132
$ws.Range("F20").Value = 1.058567155792745
$ws.Range("I20").Value = 1.043517780854394
$ws.Range("J20").Value = 1.044507889985838
$ws.Range("K20").Value = 1.046194120006132
$ws.Range("L20").Value = 1.050564414025137
$ws.Range("M20").Value = 1.061863704449394
$ws.Range("N20").Value = 1.045991210996348
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.037315505715519
$ws.Range("D21").Value = 1.042060283219968
$ws.Range("E21").Value = 1.046293979381357
$ws.Range("F21").Value = 1.057450123911388
$ws.Range("I21").Value = 1.043209569568981
$ws.Range("J21").Value = 1.043891843984333
$ws.Range("K21").Value = 1.045618012426421
$ws.Range("L21").Value = 1.049836223996125
$ws.Range("M21").Value = 1.060952220596052
$ws.Range("N21").Value = 1.045374290138859
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.036679752115428
$ws.Range("D22").Value = 1.041567650686789
$ws.Range("E22").Value = 1.045706559332746
$ws.Range("F22").Value = 1.056748957073034
$ws.Range("I22").Value = 1.043014631393071
$ws.Range("J22").Value = 1.043504292810965
$ws.Range("K22").Value = 1.045255379396951
$ws.Range("L22").Value = 1.049378550530423
$ws.Range("M22").Value = 1.060379628942105
$ws.Range("N22").Value = 1.044986188598372
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.037016731836177
$ws.Range("D23").Value = 1.041828773621815
$ws.Range("E23").Value = 1.046017879082151
$ws.Range("F23").Value = 1.057120548895162
$ws.Range("I23").Value = 1.04311808093944
$ws.Range("J23").Value = 1.043709761413936
$ws.Range("K23").Value = 1.045447656676582
$ws.Range("L23").Value = 1.049621155552259
$ws.Range("M23").Value = 1.060683123334025
$ws.Range("N23").Value = 1.045191948990327
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.03834438762036
$ws.Range("D24").Value = 1.042857467142205
$ws.Range("E24").Value = 1.047245342567573
$ws.Range("F24").Value = 1.058585902465615
$ws.Range("I24").Value = 1.043522928457241
$ws.Range("J24").Value = 1.044518214300785
$ws.Range("K24").Value = 1.046203771478288
$ws.Range("L24").Value = 1.050576625002704
$ws.Range("M24").Value = 1.061878993887797
$ws.Range("N24").Value = 1.046001549973006
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.039887083845737
$ws.Range("D25").Value = 1.044052557660458
$ws.Range("E25").Value = 1.048673476999115
$ws.Range("F25").Value = 1.060291348356869
$ws.Range("I25").Value = 1.043987687798458
$ws.Range("J25").Value = 1.045455394058194
$ws.Range("K25").Value = 1.04707937537578
$ws.Range("L25").Value = 1.051686090714523
$ws.Range("M25").Value = 1.063268846543012
$ws.Range("N25").Value = 1.046940060633167
